$d = $word.ActiveDocument

# 1) Expand the "Removed references to cileviruses" bullet into two bullets:
#    - reword the existing bullet
#    - add a new bullet (same Compact/list style) about taxonomic authorities
$bullet = $d.Content
$found = $bullet.Find.Execute(
    "Removed references to cileviruses",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We have removed all references to cileviruses^pWe have corrected inconsistencies with taxonomic authorities",
    2
)

# 2) Remove the old reviewer-comment paragraph that is no longer relevant
#    ("For the mite descriptions, there seems to be inconsistency ...", the
#    one inside the "Response to Reviewer 1" section, styled FirstParagraph)
#    and promote the following paragraph ("The wording, organization, ...")
#    to the FirstParagraph style, since it now begins this response block.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "For the mite descriptions, there seems to be inconsistency in providing taxonomic authorities*" -and $p.Style.NameLocal -eq "First Paragraph") {
        $next = $p.Next()
        $next.Style = "First Paragraph"
        $p.Range.Delete()
        break
    }
}
